$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newNote = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.37 = 13013.45 pesos`n✅ 13013.45 pesos = 3.36 = 945.42 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newNote

# --- tasas: update N10/O10 and N12/O12 rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 296.62
$ws2.Range("O10").Value = 3860.05
$ws2.Range("N12").Value = 3874.79
$ws2.Range("O12").Value = 281.5
